# Update countries & provincias Spain
#
# - Insert "Indonesia" into the country list right after "Filipinas"
#   (before "Mexico"), shifting "Mexico" and "Luxemburgo" down one row.
# - Refresh the "Datos actualizados" timestamp string.
# - Refresh the case figures for Austria (row 19), Filipinas (row 38),
#   the newly-placed Indonesia (row 39), Mexico (row 40) and
#   Luxemburgo (row 41).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp banner at the top of the sheet.
$ws.Range("A1").Value2 = "Datos actualizados a 9 de Abril de 2020 a las 10:52"

# Austria (row 19) - updated daily figures.
$ws.Range("D19").Value2 = 5240
$ws.Range("E19").Value2 = 7470
$ws.Range("F19").Value2 = 266
$ws.Range("G19").Value2 = 22
$ws.Range("H19").Value2 = 295

# Filipinas stays on row 38 but gets refreshed figures.
$ws.Range("B38").Value2 = 4076
$ws.Range("C38").Value2 = 206
$ws.Range("D38").Value2 = 124
$ws.Range("E38").Value2 = 3749
$ws.Range("G38").Value2 = 21
$ws.Range("H38").Value2 = 203

# Indonesia is inserted right after Filipinas (row 39) with new figures.
$ws.Range("A39").Value2 = "Indonesia"
$ws.Range("B39").Value2 = 3293
$ws.Range("C39").Value2 = 337
$ws.Range("D39").Value2 = 252
$ws.Range("E39").Value2 = 2761
$ws.Range("F39").Value2 = 0
$ws.Range("G39").Value2 = 40
$ws.Range("H39").Value2 = 280

# Mexico moves down to row 40, keeping its previous figures.
$ws.Range("A40").Value2 = "Mexico"
$ws.Range("B40").Value2 = 3181
$ws.Range("C40").Value2 = 396
$ws.Range("D40").Value2 = 633
$ws.Range("E40").Value2 = 2374
$ws.Range("F40").Value2 = 89
$ws.Range("G40").Value2 = 33
$ws.Range("H40").Value2 = 174

# Luxemburgo moves down to row 41, keeping its previous figures.
$ws.Range("A41").Value2 = "Luxemburgo"
$ws.Range("B41").Value2 = 3034
$ws.Range("D41").Value2 = 500
$ws.Range("E41").Value2 = 2488
$ws.Range("F41").Value2 = 34
$ws.Range("H41").Value2 = 46
